$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C12: '/' -> '20 min (fais avec detection mouvement)'
$ws.Range("C12").Value = "20 min (fais avec detection mouvement)"

# Update B13: '21 sec' -> '0:02:15.328764 (code final)'
$ws.Range("B13").Value = "0:02:15.328764 (code final)"

# Clear C13 (remove '20 min (fais avec detection mouvement)' - moved to C12)
$ws.Range("C13").ClearContents()

# Update the selection to B13 (as seen in the diff)
$ws.Range("B13").Select()
